# Apply trade-close (#84) and trade-open (#117) updates to the live
# trading results workbook: Summary, Strategy Status, All Trades and
# MarketMaking sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1400.91   # Current Capital
$summary.Range("B4").Value = 0.71      # Total P&L $
$summary.Range("B6").Value = 84        # Total Trades
$summary.Range("B8").Value = 33        # Losing Trades
$summary.Range("B9").Value = 46.43     # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.91     # Capital
$status.Range("D5").Value = 51         # Trades
$status.Range("E5").Value = 0.6        # P&L $
$status.Range("F5").Value = 0.91       # P&L %
$status.Range("G5").Value = 49.02      # Win Rate %

# ---------------------------------------------------------------------
# 3. All Trades sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# 3a. Close trade #84 (row 85)
$allTrades.Cells.Item(85, 7).Value  = 0.25              # Exit Price
$allTrades.Cells.Item(85, 8).Value  = "CLOSED"           # Status
$allTrades.Cells.Item(85, 9).Value  = -3.8462            # P&L %
$allTrades.Cells.Item(85, 10).Value = -0.01               # P&L $
$allTrades.Cells.Item(85, 11).Value = 100.91              # Capital After
$allTrades.Cells.Item(85, 12).Value = "early_exit"        # Exit Reason
$allTrades.Cells.Item(85, 13).Value = 0.15                # Duration (min)

# 3b. Append new trade #117 (row 118)
$allTrades.Cells.Item(118, 1).Value  = 117
$allTrades.Cells.Item(118, 2).Value  = "'2026-02-17"
$allTrades.Cells.Item(118, 3).Value  = "21:10:16"
$allTrades.Cells.Item(118, 4).Value  = "MarketMaking"
$allTrades.Cells.Item(118, 5).Value  = "DOWN"
$allTrades.Cells.Item(118, 6).Value  = 0.26
$allTrades.Cells.Item(118, 7).NumberFormat = "General"   # Exit Price (blank, still open)
$allTrades.Cells.Item(118, 8).Value  = "OPEN"
$allTrades.Cells.Item(118, 9).Value  = 0
$allTrades.Cells.Item(118, 10).Value = 0
$allTrades.Cells.Item(118, 11).Value = 100.9214872031006
$allTrades.Cells.Item(118, 12).NumberFormat = "General"  # Exit Reason (blank, still open)
$allTrades.Cells.Item(118, 13).Value = 0
$allTrades.Cells.Item(118, 14).Value = 0
$allTrades.Cells.Item(118, 15).Value = 0
$allTrades.Cells.Item(118, 16).Value = 0.6
$allTrades.Cells.Item(118, 17).Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------
# 4. MarketMaking sheet
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

# 4a. Close trade #84 (row 52)
$mm.Cells.Item(52, 7).Value  = 0.25                # Exit Price
$mm.Cells.Item(52, 8).Value  = "CLOSED"             # Status
$mm.Cells.Item(52, 9).Value  = -3.8462              # P&L %
$mm.Cells.Item(52, 10).Value = -0.01                 # P&L $
$mm.Cells.Item(52, 11).Value = 100.91                # Capital After
$mm.Cells.Item(52, 16).Value = "early_exit"          # Exit Reason
$mm.Cells.Item(52, 17).Value = 0.15                  # Duration (min)

# 4b. Append new trade #117 (row 85)
$mm.Cells.Item(85, 1).Value  = 117
$mm.Cells.Item(85, 2).Value  = "'2026-02-17"
$mm.Cells.Item(85, 3).Value  = "21:10:16"
$mm.Cells.Item(85, 4).Value  = "MarketMaking"
$mm.Cells.Item(85, 5).Value  = "DOWN"
$mm.Cells.Item(85, 6).Value  = 0.26
$mm.Cells.Item(85, 7).NumberFormat = "General"   # Exit Price (blank, still open)
$mm.Cells.Item(85, 8).Value  = "OPEN"
$mm.Cells.Item(85, 9).Value  = 0
$mm.Cells.Item(85, 10).Value = 0
$mm.Cells.Item(85, 11).Value = 100.9214872031006
$mm.Cells.Item(85, 12).Value = 0
$mm.Cells.Item(85, 13).Value = 0
$mm.Cells.Item(85, 14).Value = 0.6
$mm.Cells.Item(85, 15).Value = "Normal spread capture: 19600 bps"
$mm.Cells.Item(85, 16).NumberFormat = "General"  # Exit Reason (blank, still open)
$mm.Cells.Item(85, 17).Value = 0

Write-Output "Applied trade #84 close and trade #117 open updates."
